# Arbeit: Absatz 2.3.1 (Versuch) fertig
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add the "average of the three Reihen" column (H) for the first measurement
# table (rows 7-13), filled down from the formula entered in H7.
$ws.Range("H7").Formula = "=AVERAGE(D7:F7)"
$ws.Range("H7").AutoFill($ws.Range("H7:H13"), 0)

# New colour-scale conditional-formatting rule on top of the existing one.
$rng = $ws.Range("E34:G36")
$newCf = $rng.FormatConditions.AddColorScale(3)
$newCf.SetFirstPriority()
$newCf.ColorScaleCriteria.Item(1).Type = 1   # xlConditionValueLowestValue
$newCf.ColorScaleCriteria.Item(1).FormatColor.Color = 7039083   # RGB(107,105,248) -> BGR of FFF8696B
$newCf.ColorScaleCriteria.Item(2).Type = 0   # xlConditionValueNumber
$newCf.ColorScaleCriteria.Item(2).Value = -6
$newCf.ColorScaleCriteria.Item(2).FormatColor.Color = 5287058
$newCf.ColorScaleCriteria.Item(3).Type = 2   # xlConditionValueHighestValue
$newCf.ColorScaleCriteria.Item(3).FormatColor.Color = 39423

# Select H7:H13 with the active cell on H7, matching the final selection
# state left behind after filling the formula down.
$ws.Range("H7:H13").Select()
$ws.Application.ActiveWindow.RangeSelection.Item(1).Activate()
